$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 34800
$ws.Cells.Item(3, 10).Value = 34800
$ws.Cells.Item(3, 12).Value = 34800
$ws.Cells.Item(3, 14).Value = -35028
$ws.Cells.Item(51, 8).Value = 6870
$ws.Cells.Item(51, 9).Value = 2000
$ws.Cells.Item(51, 10).Value = 8493.333000000001
$ws.Cells.Item(51, 11).Value = 2000
$ws.Cells.Item(51, 12).Value = 8493.333000000001
$ws.Cells.Item(51, 13).Value = -1516
$ws.Cells.Item(51, 14).Value = -9461.333000000001
$ws.Cells.Item(80, 8).Value = 3082756.8
$ws.Cells.Item(80, 9).Value = 1780.6
$ws.Cells.Item(80, 10).Value = 3853000.8
$ws.Cells.Item(80, 11).Value = 5341.799999999999
$ws.Cells.Item(80, 12).Value = 11559002.4
$ws.Cells.Item(80, 13).Value = -4343.799999999999
$ws.Cells.Item(80, 14).Value = -11560998.4
$ws.Cells.Item(83, 8).Value = 3082756.8
$ws.Cells.Item(83, 9).Value = 1780.6
$ws.Cells.Item(83, 10).Value = 3853000.8
$ws.Cells.Item(83, 11).Value = 16025.4
$ws.Cells.Item(83, 12).Value = 34677007.2
$ws.Cells.Item(83, 13).Value = -11033.4
$ws.Cells.Item(83, 14).Value = -34686991.2
$ws.Cells.Item(94, 8).Value = 0
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).ClearContents()
$ws.Cells.Item(94, 14).ClearContents()
$ws.Cells.Item(97, 8).Value = 1025
$ws.Cells.Item(97, 10).Value = 1025
$ws.Cells.Item(97, 12).Value = 3075
$ws.Cells.Item(97, 14).Value = -4067
$ws.Cells.Item(98, 8).Value = 658.7857
$ws.Cells.Item(98, 9).Value = 626.4
$ws.Cells.Item(98, 10).Value = 739.75
$ws.Cells.Item(98, 11).Value = 626.4
$ws.Cells.Item(98, 12).Value = 739.75
$ws.Cells.Item(98, 13).Value = 871.6
$ws.Cells.Item(98, 14).Value = -3735.75
$ws.Cells.Item(100, 8).Value = 1663.6
$ws.Cells.Item(100, 9).Value = 1520
$ws.Cells.Item(100, 10).Value = 1998.6666
$ws.Cells.Item(100, 11).Value = 1520
$ws.Cells.Item(100, 12).Value = 1998.6666
$ws.Cells.Item(100, 13).Value = -979
$ws.Cells.Item(100, 14).Value = -3080.6666
$ws.Cells.Item(102, 8).Value = 34800
$ws.Cells.Item(102, 10).Value = 34800
$ws.Cells.Item(102, 12).Value = 34800
$ws.Cells.Item(102, 14).Value = -41290
$ws.Cells.Item(103, 8).Value = 139152.78
$ws.Cells.Item(103, 9).Value = 156497.5
$ws.Cells.Item(103, 10).Value = 395
$ws.Cells.Item(103, 11).Value = 469492.5
$ws.Cells.Item(103, 12).Value = 1185
$ws.Cells.Item(103, 13).Value = -468906.5
$ws.Cells.Item(103, 14).Value = -2357
$ws.Cells.Item(122, 8).Value = 658.7857
$ws.Cells.Item(122, 9).Value = 626.4
$ws.Cells.Item(122, 10).Value = 739.75
$ws.Cells.Item(122, 11).Value = 1879.2
$ws.Cells.Item(122, 12).Value = 2219.25
$ws.Cells.Item(122, 13).Value = 570.8000000000002
$ws.Cells.Item(122, 14).Value = -7119.25
$ws.Cells.Item(132, 8).Value = 3350.76
$ws.Cells.Item(132, 9).Value = 3773.75
$ws.Cells.Item(132, 10).Value = 1658.8
$ws.Cells.Item(132, 11).Value = 11321.25
$ws.Cells.Item(132, 12).Value = 4976.4
$ws.Cells.Item(132, 13).Value = -8791.25
$ws.Cells.Item(132, 14).Value = -10036.4
$ws.Cells.Item(137, 8).Value = 32358.363
$ws.Cells.Item(137, 9).Value = 2125.75
$ws.Cells.Item(137, 10).Value = 78870.08
$ws.Cells.Item(137, 11).Value = 6377.25
$ws.Cells.Item(137, 12).Value = 236610.24
$ws.Cells.Item(137, 13).Value = -3827.25
$ws.Cells.Item(137, 14).Value = -241710.24

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2356.1428
$ws.Cells.Item(2, 9).Value = 1442.3846
$ws.Cells.Item(2, 10).Value = 3841
$ws.Cells.Item(2, 11).Value = 1442.3846
$ws.Cells.Item(2, 12).Value = 3841
$ws.Cells.Item(2, 13).Value = -1329.3846
$ws.Cells.Item(2, 14).Value = -4067
$ws.Cells.Item(32, 8).Value = 19706.197
$ws.Cells.Item(32, 9).Value = 21362
$ws.Cells.Item(32, 10).Value = 2817
$ws.Cells.Item(32, 11).Value = 21362
$ws.Cells.Item(32, 12).Value = 2817
$ws.Cells.Item(32, 13).Value = -21075
$ws.Cells.Item(32, 14).Value = -3391
$ws.Cells.Item(45, 8).Value = 3008.3914
$ws.Cells.Item(45, 10).Value = 3054.0356
$ws.Cells.Item(45, 12).Value = 3054.0356
$ws.Cells.Item(45, 14).Value = -3808.0356
$ws.Cells.Item(80, 8).Value = 47853.75
$ws.Cells.Item(80, 10).Value = 47853.75
$ws.Cells.Item(80, 12).Value = 47853.75
$ws.Cells.Item(80, 14).Value = -49849.75
$ws.Cells.Item(83, 8).Value = 47853.75
$ws.Cells.Item(83, 10).Value = 47853.75
$ws.Cells.Item(83, 12).Value = 143561.25
$ws.Cells.Item(83, 14).Value = -153545.25
$ws.Cells.Item(88, 8).Value = 75080.78999999999
$ws.Cells.Item(88, 10).Value = 95011.91
$ws.Cells.Item(88, 12).Value = 95011.91
$ws.Cells.Item(88, 14).Value = -95823.91
$ws.Cells.Item(91, 8).Value = 75080.78999999999
$ws.Cells.Item(91, 10).Value = 95011.91
$ws.Cells.Item(91, 12).Value = 95011.91
$ws.Cells.Item(91, 14).Value = -97819.91
$ws.Cells.Item(110, 8).Value = 2192.4443
$ws.Cells.Item(110, 9).Value = 0
$ws.Cells.Item(110, 10).Value = 2192.4443
$ws.Cells.Item(110, 11).Value = 0
$ws.Cells.Item(110, 12).Value = 2192.4443
$ws.Cells.Item(110, 13).ClearContents()
$ws.Cells.Item(110, 14).Value = -6282.4443
$ws.Cells.Item(116, 8).Value = 2356.1428
$ws.Cells.Item(116, 9).Value = 1442.3846
$ws.Cells.Item(116, 10).Value = 3841
$ws.Cells.Item(116, 11).Value = 1442.3846
$ws.Cells.Item(116, 12).Value = 3841
$ws.Cells.Item(116, 13).Value = 851.6153999999999
$ws.Cells.Item(116, 14).Value = -8429
$ws.Cells.Item(122, 8).Value = 1890.2222
$ws.Cells.Item(122, 9).Value = 1939
$ws.Cells.Item(122, 10).Value = 1500
$ws.Cells.Item(122, 11).Value = 5817
$ws.Cells.Item(122, 12).Value = 4500
$ws.Cells.Item(122, 13).Value = -3367
$ws.Cells.Item(122, 14).Value = -9400

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2356.1428
$ws.Cells.Item(3, 9).Value = 1442.3846
$ws.Cells.Item(3, 10).Value = 3841
$ws.Cells.Item(3, 11).Value = 1442.3846
$ws.Cells.Item(3, 12).Value = 3841
$ws.Cells.Item(3, 13).Value = -1328.3846
$ws.Cells.Item(3, 14).Value = -4069
$ws.Cells.Item(86, 8).Value = 2422.7856
$ws.Cells.Item(86, 9).Value = 1740.5
$ws.Cells.Item(86, 10).Value = 4128.5
$ws.Cells.Item(86, 11).Value = 1740.5
$ws.Cells.Item(86, 12).Value = 4128.5
$ws.Cells.Item(86, 13).Value = -617.5
$ws.Cells.Item(86, 14).Value = -6374.5
$ws.Cells.Item(89, 8).Value = 2422.7856
$ws.Cells.Item(89, 9).Value = 1740.5
$ws.Cells.Item(89, 10).Value = 4128.5
$ws.Cells.Item(89, 11).Value = 8702.5
$ws.Cells.Item(89, 12).Value = 20642.5
$ws.Cells.Item(89, 13).Value = -3086.5
$ws.Cells.Item(89, 14).Value = -31874.5
$ws.Cells.Item(99, 8).Value = 2524.1428
$ws.Cells.Item(99, 9).Value = 2114.8333
$ws.Cells.Item(99, 10).Value = 4980
$ws.Cells.Item(99, 11).Value = 2114.8333
$ws.Cells.Item(99, 12).Value = 4980
$ws.Cells.Item(99, 13).Value = -616.8332999999998
$ws.Cells.Item(99, 14).Value = -7976
$ws.Cells.Item(134, 8).Value = 45208.207
$ws.Cells.Item(134, 9).Value = 53949.85
$ws.Cells.Item(134, 10).Value = 1500
$ws.Cells.Item(134, 11).Value = 161849.55
$ws.Cells.Item(134, 12).Value = 4500
$ws.Cells.Item(134, 13).Value = -159314.55
$ws.Cells.Item(134, 14).Value = -9570

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 6013.1113
$ws.Cells.Item(62, 9).Value = 4275
$ws.Cells.Item(62, 11).Value = 4275
$ws.Cells.Item(62, 13).Value = -3651
$ws.Cells.Item(65, 8).Value = 6013.1113
$ws.Cells.Item(65, 9).Value = 4275
$ws.Cells.Item(65, 11).Value = 21375
$ws.Cells.Item(65, 13).Value = -18255
$ws.Cells.Item(86, 8).Value = 17157.3
$ws.Cells.Item(86, 9).Value = 3226.75
$ws.Cells.Item(86, 10).Value = 26444.334
$ws.Cells.Item(86, 11).Value = 3226.75
$ws.Cells.Item(86, 12).Value = 26444.334
$ws.Cells.Item(86, 13).Value = -2103.75
$ws.Cells.Item(86, 14).Value = -28690.334
$ws.Cells.Item(89, 8).Value = 17157.3
$ws.Cells.Item(89, 9).Value = 3226.75
$ws.Cells.Item(89, 10).Value = 26444.334
$ws.Cells.Item(89, 11).Value = 16133.75
$ws.Cells.Item(89, 12).Value = 132221.67
$ws.Cells.Item(89, 13).Value = -10517.75
$ws.Cells.Item(89, 14).Value = -143453.67
$ws.Cells.Item(134, 8).Value = 999.9697
$ws.Cells.Item(134, 9).Value = 761.0454999999999
$ws.Cells.Item(134, 10).Value = 1477.8182
$ws.Cells.Item(134, 11).Value = 2283.1365
$ws.Cells.Item(134, 12).Value = 4433.4546
$ws.Cells.Item(134, 13).Value = 251.8635000000004
$ws.Cells.Item(134, 14).Value = -9503.454600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(9, 8).Value = 2101
$ws.Cells.Item(9, 9).Value = 301
$ws.Cells.Item(9, 10).Value = 3001
$ws.Cells.Item(9, 11).Value = 903
$ws.Cells.Item(9, 12).Value = 9003
$ws.Cells.Item(9, 13).Value = -679
$ws.Cells.Item(9, 14).Value = -9451
$ws.Cells.Item(26, 8).Value = 509
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 509
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 1527
$ws.Cells.Item(26, 13).ClearContents()
$ws.Cells.Item(26, 14).Value = -2103
$ws.Cells.Item(76, 8).Value = 4473.846
$ws.Cells.Item(76, 9).Value = 1500
$ws.Cells.Item(76, 10).Value = 5014.5454
$ws.Cells.Item(76, 11).Value = 4500
$ws.Cells.Item(76, 12).Value = 15043.6362
$ws.Cells.Item(76, 13).Value = -4117
$ws.Cells.Item(76, 14).Value = -15809.6362
$ws.Cells.Item(79, 8).Value = 4473.846
$ws.Cells.Item(79, 9).Value = 1500
$ws.Cells.Item(79, 10).Value = 5014.5454
$ws.Cells.Item(79, 11).Value = 4500
$ws.Cells.Item(79, 12).Value = 15043.6362
$ws.Cells.Item(79, 13).Value = -3174
$ws.Cells.Item(79, 14).Value = -17695.6362
$ws.Cells.Item(98, 8).Value = 900
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 900
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 12).Value = 2700
$ws.Cells.Item(98, 13).ClearContents()
$ws.Cells.Item(98, 14).Value = -5696
$ws.Cells.Item(104, 8).Value = 10029
$ws.Cells.Item(104, 10).Value = 10029
$ws.Cells.Item(104, 12).Value = 30087
$ws.Cells.Item(104, 14).Value = -35329
$ws.Cells.Item(123, 8).Value = 3881.8333
$ws.Cells.Item(123, 9).Value = 1001
$ws.Cells.Item(123, 10).Value = 4458
$ws.Cells.Item(123, 11).Value = 3003
$ws.Cells.Item(123, 12).Value = 13374
$ws.Cells.Item(123, 13).Value = -553
$ws.Cells.Item(123, 14).Value = -18274
$ws.Cells.Item(131, 8).Value = 747.96
$ws.Cells.Item(131, 10).Value = 763.8125
$ws.Cells.Item(131, 12).Value = 2291.4375
$ws.Cells.Item(131, 14).Value = -12371.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1772.3572
$ws.Cells.Item(97, 9).Value = 1208.55
$ws.Cells.Item(97, 11).Value = 1208.55
$ws.Cells.Item(97, 13).Value = -712.55
$ws.Cells.Item(102, 8).Value = 1746.5217
$ws.Cells.Item(102, 9).Value = 1790.8823
$ws.Cells.Item(102, 10).Value = 1620.8334
$ws.Cells.Item(102, 11).Value = 1790.8823
$ws.Cells.Item(102, 12).Value = 1620.8334
$ws.Cells.Item(102, 13).Value = -168.8823
$ws.Cells.Item(102, 14).Value = -4864.8334
$ws.Cells.Item(113, 8).Value = 3123.75
$ws.Cells.Item(113, 9).Value = 1978
$ws.Cells.Item(113, 11).Value = 1978
$ws.Cells.Item(113, 13).Value = 192
$ws.Cells.Item(122, 8).Value = 1623.4615
$ws.Cells.Item(122, 9).Value = 1585.5714
$ws.Cells.Item(122, 10).Value = 1667.6666
$ws.Cells.Item(122, 11).Value = 4756.7142
$ws.Cells.Item(122, 12).Value = 5002.9998
$ws.Cells.Item(122, 13).Value = -2306.7142
$ws.Cells.Item(122, 14).Value = -9902.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 748.13043
$ws.Cells.Item(100, 9).Value = 412
$ws.Cells.Item(100, 10).Value = 1700.5
$ws.Cells.Item(100, 11).Value = 824
$ws.Cells.Item(100, 12).Value = 3401
$ws.Cells.Item(100, 13).Value = -283
$ws.Cells.Item(100, 14).Value = -4483
$ws.Cells.Item(126, 8).Value = 2461
$ws.Cells.Item(126, 10).Value = 2641.1428
$ws.Cells.Item(126, 12).Value = 7923.428400000001
$ws.Cells.Item(126, 14).Value = -12863.4284
